$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F2").Value = 61
$ws.Range("H2").Value = 61
$ws.Range("F5").Value = 80
$ws.Range("H5").Value = 80
$ws.Range("F6").Value = 27
$ws.Range("H6").Value = 27
$ws.Range("E10").Value = 413
$ws.Range("F10").Value = 206
$ws.Range("H10").Value = 206
$ws.Range("F11").Value = 152
$ws.Range("H11").Value = 152
$ws.Range("E12").Value = 402
$ws.Range("F12").Value = 220
$ws.Range("H12").Value = 220
$ws.Range("E13").Value = 106
$ws.Range("E15").Value = 133
$ws.Range("F15").Value = 53
$ws.Range("H15").Value = 53
$ws.Range("E16").Value = 170
$ws.Range("F16").Value = 88
$ws.Range("H16").Value = 88
$ws.Range("E17").Value = 77
$ws.Range("F17").Value = 37
$ws.Range("H17").Value = 37
$ws.Range("F18").Value = 22
$ws.Range("H18").Value = 22
$ws.Range("F20").Value = 28
$ws.Range("H20").Value = 28
$ws.Range("F22").Value = 74
$ws.Range("H22").Value = 74
$ws.Range("E23").Value = 173
$ws.Range("E25").Value = 211
$ws.Range("F25").Value = 99
$ws.Range("H25").Value = 99
$ws.Range("E26").Value = 119
$ws.Range("F26").Value = 75
$ws.Range("H26").Value = 75
$ws.Range("E27").Value = 273
$ws.Range("F27").Value = 129
$ws.Range("H27").Value = 129
$ws.Range("F29").Value = 80
$ws.Range("H29").Value = 80
$ws.Range("F30").Value = 98
$ws.Range("H30").Value = 98
$ws.Range("E31").Value = 66
$ws.Range("F31").Value = 30
$ws.Range("H31").Value = 30
$ws.Range("F33").Value = 125
$ws.Range("H33").Value = 125
$ws.Range("F34").Value = 110
$ws.Range("H34").Value = 110
$ws.Range("F35").Value = 76
$ws.Range("H35").Value = 76
$ws.Range("E36").Value = 54
$ws.Range("E37").Value = 133
$ws.Range("F37").Value = 63
$ws.Range("H37").Value = 63
$ws.Range("E38").Value = 82
$ws.Range("E39").Value = 163
$ws.Range("F39").Value = 77
$ws.Range("H39").Value = 77
$ws.Range("F40").Value = 102
$ws.Range("H40").Value = 102
$ws.Range("F41").Value = 154
$ws.Range("H41").Value = 154
$ws.Range("E42").Value = 301
$ws.Range("F42").Value = 160
$ws.Range("H42").Value = 160
$ws.Range("E44").Value = 257
$ws.Range("F44").Value = 126
$ws.Range("H44").Value = 126
$ws.Range("F45").Value = 57
$ws.Range("H45").Value = 57
$ws.Range("E46").Value = 263
$ws.Range("F46").Value = 144
$ws.Range("H46").Value = 144
$ws.Range("E47").Value = 372
$ws.Range("F47").Value = 184
$ws.Range("H47").Value = 184
$ws.Range("F48").Value = 70
$ws.Range("H48").Value = 70
$ws.Range("E49").Value = 251
$ws.Range("E50").Value = 216
$ws.Range("F50").Value = 93
$ws.Range("H50").Value = 93
$ws.Range("E51").Value = 202
$ws.Range("F51").Value = 81
$ws.Range("H51").Value = 81